$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: remove "Bakda Pratiwi", fix membership-number typo (UL -> Ul),
#     and change the tagihan (bill) computation method ---
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "2083/Kopma_Ul/21"
$ws.Range("D2").Value = "Anggota"
$ws.Range("H2").Value = 10000

# --- Row 3: Budi Gunawan - tagihan method changed ---
$ws.Range("F3").Value = 30000
$ws.Range("G3").Value = 60000
$ws.Range("H3").Value = 5000

# --- Insert a new row 4 for the new member "Rian" (everything below shifts
#     down by one row, so the last member becomes the new Alumni entry) ---
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Rian"
$ws.Range("C4").Value = "2086/Kopma_UL/20"
$ws.Range("D4").Value = "-- Status Keanggotaan --"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

# --- Row 5 (was row 4): Hartana Prima Prayoga S.Kom - tagihan method changed ---
$ws.Range("A5").Value = 4
$ws.Range("F5").Value = 150000
$ws.Range("G5").Value = 150000
$ws.Range("H5").Value = 0

# --- Row 6 (was row 5): Niyaga Suryono - tagihan method changed ---
$ws.Range("A6").Value = 5
$ws.Range("F6").Value = 50000
$ws.Range("G6").Value = 50000
$ws.Range("H6").Value = 0

# --- Rows 7-12 (were rows 6-11): remaining members, tagihan zeroed out
#     under the new method, and renumber the "No" column ---
$ws.Range("A7").Value = 6
$ws.Range("H7").Value = 0

$ws.Range("A8").Value = 7
$ws.Range("H8").Value = 0

$ws.Range("A9").Value = 8
$ws.Range("H9").Value = 0

$ws.Range("A10").Value = 9
$ws.Range("H10").Value = 0

$ws.Range("A11").Value = 10
$ws.Range("H11").Value = 0

$ws.Range("A12").Value = 11
$ws.Range("H12").Value = 0
